$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 15151898
$ws.Range("I12").Value = 45454544
$ws.Range("K12").Value = 45454544
$ws.Range("M12").Value = -45454374
$ws.Range("H129").Value = 1464.2727
$ws.Range("J129").Value = 2600
$ws.Range("L129").Value = 7800
$ws.Range("N129").Value = -17800
$ws.Range("H137").Value = 3134460.5
$ws.Range("I137").Value = 7965.125
$ws.Range("K137").Value = 23895.375
$ws.Range("M137").Value = -21345.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1715978.5
$ws.Range("I61").Value = 69885.56
$ws.Range("K61").Value = 69885.56
$ws.Range("M61").Value = -69673.56
$ws.Range("H74").Value = 637599.1
$ws.Range("I74").Value = 1998.4482
$ws.Range("J74").Value = 3709669.2
$ws.Range("K74").Value = 1998.4482
$ws.Range("L74").Value = 3709669.2
$ws.Range("M74").Value = -1124.4482
$ws.Range("N74").Value = -3711417.2
$ws.Range("H77").Value = 637599.1
$ws.Range("I77").Value = 1998.4482
$ws.Range("J77").Value = 3709669.2
$ws.Range("K77").Value = 9992.241
$ws.Range("L77").Value = 18548346
$ws.Range("M77").Value = -5624.241
$ws.Range("N77").Value = -18557082
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 2920.7368
$ws.Range("I132").Value = 2716.4285
$ws.Range("J132").Value = 3039.9167
$ws.Range("K132").Value = 8149.2855
$ws.Range("L132").Value = 9119.750100000001
$ws.Range("M132").Value = -5619.2855
$ws.Range("N132").Value = -14179.7501
$ws.Range("H136").Value = 1715978.5
$ws.Range("I136").Value = 69885.56
$ws.Range("K136").Value = 209656.68
$ws.Range("M136").Value = -207106.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H104").Value = 60000
$ws.Range("J104").Value = 60000
$ws.Range("L104").Value = 60000
$ws.Range("N104").Value = -66988
$ws.Range("H105").Value = 8779.963
$ws.Range("I105").Value = 7397.3335
$ws.Range("J105").Value = 11545.223
$ws.Range("K105").Value = 7397.3335
$ws.Range("L105").Value = 11545.223
$ws.Range("M105").Value = -5650.3335
$ws.Range("N105").Value = -15039.223
$ws.Range("H134").Value = 28127006
$ws.Range("I134").Value = 2251.8572
$ws.Range("K134").Value = 6755.571599999999
$ws.Range("M134").Value = -4220.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 590.7778
$ws.Range("J22").Value = 877.6667
$ws.Range("L22").Value = 877.6667
$ws.Range("N22").Value = -1577.6667
$ws.Range("H31").Value = 6544.905
$ws.Range("I31").Value = 5465.5
$ws.Range("J31").Value = 9999
$ws.Range("K31").Value = 5465.5
$ws.Range("L31").Value = 9999
$ws.Range("M31").Value = -5170.5
$ws.Range("N31").Value = -10589
$ws.Range("H34").Value = 6544.905
$ws.Range("I34").Value = 5465.5
$ws.Range("J34").Value = 9999
$ws.Range("K34").Value = 5465.5
$ws.Range("L34").Value = 9999
$ws.Range("M34").Value = -5263.5
$ws.Range("N34").Value = -10403
$ws.Range("H58").Value = 1899.4706
$ws.Range("J58").Value = 2375.5557
$ws.Range("L58").Value = 2375.5557
$ws.Range("N58").Value = -2781.5557
$ws.Range("H62").Value = 4202.9165
$ws.Range("J62").Value = 5828
$ws.Range("L62").Value = 5828
$ws.Range("N62").Value = -7076
$ws.Range("H65").Value = 4202.9165
$ws.Range("J65").Value = 5828
$ws.Range("L65").Value = 29140
$ws.Range("N65").Value = -35380
$ws.Range("H86").Value = 9365.261
$ws.Range("I86").Value = 4088
$ws.Range("K86").Value = 4088
$ws.Range("M86").Value = -2965
$ws.Range("H89").Value = 9365.261
$ws.Range("I89").Value = 4088
$ws.Range("K89").Value = 20440
$ws.Range("M89").Value = -14824
$ws.Range("H107").Value = 1379.64
$ws.Range("I107").Value = 1269.7
$ws.Range("J107").Value = 1819.4
$ws.Range("K107").Value = 1269.7
$ws.Range("L107").Value = 1819.4
$ws.Range("M107").Value = 650.3
$ws.Range("N107").Value = -5659.4
$ws.Range("H136").Value = 1899.4706
$ws.Range("J136").Value = 2375.5557
$ws.Range("L136").Value = 7126.6671
$ws.Range("N136").Value = -12226.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 718851.5600000001
$ws.Range("I4").Value = 272.80487
$ws.Range("K4").Value = 818.41461
$ws.Range("M4").Value = -706.41461
$ws.Range("H17").Value = 1200
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1200
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3600
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3938
$ws.Range("H103").Value = 1114836.9
$ws.Range("J103").Value = 10833
$ws.Range("L103").Value = 32499
$ws.Range("N103").Value = -34257
$ws.Range("H107").Value = 674.575
$ws.Range("I107").Value = 589.43475
$ws.Range("K107").Value = 1768.30425
$ws.Range("M107").Value = 151.6957499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8153603
$ws.Range("I80").Value = 99399.35000000001
$ws.Range("K80").Value = 99399.35000000001
$ws.Range("M80").Value = -98401.35000000001
$ws.Range("H83").Value = 8153603
$ws.Range("I83").Value = 99399.35000000001
$ws.Range("K83").Value = 496996.75
$ws.Range("M83").Value = -492004.75
$ws.Range("H113").Value = 4766.375
$ws.Range("I113").Value = 4545.1816
$ws.Range("K113").Value = 4545.1816
$ws.Range("M113").Value = -2375.1816
$ws.Range("H122").Value = 3271.0293
$ws.Range("I122").Value = 2332.6072
$ws.Range("J122").Value = 7650.3335
$ws.Range("K122").Value = 6997.821599999999
$ws.Range("L122").Value = 22951.0005
$ws.Range("M122").Value = -4547.821599999999
$ws.Range("N122").Value = -27851.0005
$ws.Range("H126").Value = 4754.0835
$ws.Range("I126").Value = 4209.125
$ws.Range("K126").Value = 12627.375
$ws.Range("M126").Value = -10157.375
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 4292995.5
$ws.Range("I132").Value = 23390
$ws.Range("J132").Value = 25641024
$ws.Range("K132").Value = 70170
$ws.Range("L132").Value = 76923072
$ws.Range("M132").Value = -67640
$ws.Range("N132").Value = -76928132
$ws.Range("H134").Value = 94000
$ws.Range("J134").Value = 94000
$ws.Range("L134").Value = 282000
$ws.Range("N134").Value = -287070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7520.091
$ws.Range("I46").Value = 22795.4
$ws.Range("K46").Value = 22795.4
$ws.Range("M46").Value = -22607.4
$ws.Range("H82").Value = 813.125
$ws.Range("I82").Value = 793.2308
$ws.Range("K82").Value = 793.2308
$ws.Range("M82").Value = -432.2308
$ws.Range("H85").Value = 813.125
$ws.Range("I85").Value = 793.2308
$ws.Range("K85").Value = 793.2308
$ws.Range("M85").Value = 454.7692
$ws.Range("H93").Value = 1329.9166
$ws.Range("I93").Value = 1462.3334
$ws.Range("J93").Value = 932.6667
$ws.Range("K93").Value = 1462.3334
$ws.Range("L93").Value = 932.6667
$ws.Range("M93").Value = -214.3334
$ws.Range("N93").Value = -3428.6667
$ws.Range("H100").Value = 3072.5
$ws.Range("I100").Value = 2715.4546
$ws.Range("K100").Value = 2715.4546
$ws.Range("M100").Value = -2174.4546
$ws.Range("H132").Value = 3370.4
$ws.Range("I132").Value = 957.3333
$ws.Range("K132").Value = 2871.9999
$ws.Range("M132").Value = -341.9998999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 7663
$ws.Range("J23").Value = 10994.5
$ws.Range("L23").Value = 10994.5
$ws.Range("N23").Value = -11452.5
$ws.Range("H107").Value = 773906.5600000001
$ws.Range("J107").Value = 1589053.1
$ws.Range("L107").Value = 4767159.300000001
$ws.Range("N107").Value = -4770999.300000001
$ws.Range("H132").Value = 1558.2778
$ws.Range("I132").Value = 1408.7561
$ws.Range("K132").Value = 4226.2683
$ws.Range("M132").Value = -1696.2683
